$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1 (paragraph "Every time we make a decision..."):
#   The source diff wraps the existing "make a decision" run with
#   <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/>
#   markers. Those are Word's internal grammar-checker breadcrumbs; the
#   run/text itself is unchanged, so there is nothing to edit here beyond
#   leaving the run boundary intact (it already is its own run).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Edit 2 (same paragraph): "those alternatives have tradeoffs that need to
# be considered." -> "those alternatives and tradeoffs need to be
# considered.", with the result split across three runs:
#   "hose alternatives " | "and" | " tradeoffs need to be considered."
# ---------------------------------------------------------------------------

# "have" -> "and"
$find1 = $d.Content
$find1.Find.Execute("have tradeoffs")
$haveRange = $d.Range($find1.Start, $find1.Start + 4)
$haveRange.Font.Bold = 1
$haveRange.Text = "and"
$haveRange.Font.Bold = 0

# remove " that" from " tradeoffs that need"
$find2 = $d.Content
$find2.Find.Execute("tradeoffs that need")
$thatStart = $find2.Start + 10
$thatRange = $d.Range($thatStart, $thatStart + 5)
$thatRange.Font.Bold = 1
$thatRange.Text = ""
$thatRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 3 (paragraph "Ethics is why I dove into..."): split the run ending in
# "...Aquinas actually provide" so "actually provide" becomes its own run
# (it is immediately followed by the existing separate "d" run). In the
# source diff both runs end up flanked by gramStart/gramEnd proofErr marks.
# ---------------------------------------------------------------------------
$find3 = $d.Content
$find3.Find.Execute("actually provide")
$apRange = $d.Range($find3.Start, $find3.End)
$apRange.Font.Bold = 1
$apRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 4 (same paragraph): "...basis, we actually can. " gets split into
#   "...basis, we " | "actually can" | ". "
# ---------------------------------------------------------------------------
$find4 = $d.Content
$find4.Find.Execute("actually can")
$acRange = $d.Range($find4.Start, $find4.End)
$acRange.Font.Bold = 1
$acRange.Font.Bold = 0
